$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two data rows (2 and 3) had their weekly date and price values swapped.
# Row 2 becomes what used to be row 3's values, and vice versa, for columns
# D (Fecha), M (Volumen), N (Precio mínimo), O (Precio máximo),
# P (Precio promedio ponderado) and S (Precio $/Kg).

$ws.Range("D2").Value = 44365
$ws.Range("M2").Value = 900
$ws.Range("N2").Value = 1200
$ws.Range("O2").Value = 1400
$ws.Range("P2").Value = 1300
$ws.Range("S2").Value = 1300

$ws.Range("D3").Value = 44322
$ws.Range("M3").Value = 600
$ws.Range("N3").Value = 1500
$ws.Range("O3").Value = 1600
$ws.Range("P3").Value = 1550
$ws.Range("S3").Value = 1550
